$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.654.46'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.597.44'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '211.39'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.822.29'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '1.544.49'
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").Value = '0.523'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '65.10'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '26.649.39'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("E18").Value = '  +1.32%  '
$ws.Range("D19").Value = '209.22'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = '145.07'
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '7.10'
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("D29").Value = '15.28'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").Value = '0.0514'
$ws.Range("E30").Value = '  +1.87%  '
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").Value = '1.284.73'
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").Value = '0.618'
$ws.Range("E35").Value = '  -7.14%  '
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("D39").Value = '1.07'
$ws.Range("E39").Value = '  +22.21%  '
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").Value = '0.784'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '63.88'
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").Value = '1.734.60'
$ws.Range("D46").Value = '90.79'
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("E48").Value = '  +2.46%  '
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = '7.41'
$ws.Range("E51").Value = '  -1.24%  '
